$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "0.015±0.003"
$ws.Range("C2").Value = "0.199±0.008"

$ws.Range("B3").Value = "0.097±0.003"
$ws.Range("C3").Value = "0.225±0.045"

$ws.Range("B4").Value = "0.530±0.104"
$ws.Range("C4").Value = "0.124±0.021"

$ws.Range("B5").Value = "0.959±0.011"
$ws.Range("C5").Value = "0.410±0.027"

$ws.Range("B6").Value = "0.785±0.050"
$ws.Range("C6").Value = "0.539±0.084"

$ws.Range("B7").Value = "0.561±0.104"
$ws.Range("C7").Value = "0.069±0.022"

$ws.Range("B8").Value = "0.007±0.002"
$ws.Range("C8").Value = "0.242±0.020"

$ws.Range("B9").Value = "0.077±0.004"
$ws.Range("C9").Value = "0.172±0.032"

$ws.Range("B10").Value = "0.781±0.045"
$ws.Range("C10").Value = "0.571±0.062"
